$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells: AD1=Wins, AE1=Losses, AF1=Ties ---
$ws.Range("AD1").Value2 = "Wins"
$ws.Range("AE1").Value2 = "Losses"
$ws.Range("AF1").Value2 = "Ties"

# Match the formatting used by the rest of the header row (bold, bordered,
# centered) by copying A1's format onto the three new header cells.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Team record columns for every data row (2-46): Wins=96, Losses=66, Ties=0 ---
$lastRow = 46
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value2 = 96
    $ws.Cells.Item($r, 31).Value2 = 66
    $ws.Cells.Item($r, 32).Value2 = 0
}
